{"js": "// The document contains a paragraph built from three separate runs whose\n// text concatenates to \"<id>p108r_2</id>\" (the runs are \"<id>\", \"p108r_2\",\n// \"</id>\"). The edit merges these three runs into a single run containing\n// the full literal text \"<id>p108r_2</id>\", keeping the formatting of the\n// first (\"<id>\") run (Courier New / color 7f6000 / size 9pt) and leaving\n// the rest of the paragraph (including the trailing empty run) untouched.\n\nconst body = context.document.body;\n\n// Locate the text by searching across run boundaries for the exact\n// (pre-edit) concatenation of the three runs.\nconst results = body.search(\"<id>p108r_2</id>\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"<id>p108r_2</id>\" in the document body.');\n}\n\n// Replacing the matched range's text collapses it into a single run that\n// inherits the formatting of the range's leading run (the \"<id>\" run),\n// exactly matching the target OOXML.\nconst target = results.items[0];\ntarget.insertText(\"<id>p108r_2</id>\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The document contains a paragraph built from three separate runs whose\n# text concatenates to \"<id>p108r_2</id>\" (the runs are \"<id>\", \"p108r_2\",\n# \"</id>\"). The edit merges these three runs into a single run containing\n# the full literal text \"<id>p108r_2</id>\", keeping the formatting of the\n# first (\"<id>\") run (Courier New / color 7f6000 / size 9pt) and leaving\n# the rest of the paragraph (including the trailing empty run) untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the trailing two runs' text (\"p108r_2\" + \"</id>\"), which sit\n#    immediately after the \"<id>\" run. Deleting this span merges it away,\n#    leaving only the \"<id>\" run (with its original formatting/rsids).\n$tail = $d.Content\n$foundTail = $tail.Find.Execute(\"p108r_2</id>\")\nif (-not $foundTail) {\n    throw \"Could not find 'p108r_2</id>' in the document.\"\n}\n$tail.Delete()\n\n# 2) Re-locate the now-standalone \"<id>\" run and append the removed text\n#    back onto it, so it becomes a single run containing the full\n#    \"<id>p108r_2</id>\" string with the \"<id>\" run's formatting.\n$head = $d.Content\n$foundHead = $head.Find.Execute(\"<id>\")\nif (-not $foundHead) {\n    throw \"Could not find '<id>' in the document.\"\n}\n$head.Collapse(0)  # wdCollapseEnd\n$head.InsertAfter(\"p108r_2</id>\")\n"}
